$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.785.91'
$ws.Range("E2").Value = '  -3.72%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.023.91'
$ws.Range("E3").Value = '  -3.02%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.91'
$ws.Range("E5").Value = '  +0.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.46'
$ws.Range("E6").Value = '  -3.12%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.020.73'
$ws.Range("E8").Value = '  -2.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.497'
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.149'
$ws.Range("E10").Value = '  -5.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.07'
$ws.Range("E11").Value = '  -6.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  -2.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000222'
$ws.Range("E13").Value = '  -2.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.38'
$ws.Range("E14").Value = '  -2.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.518.42'
$ws.Range("E15").Value = '  -2.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '61.980.41'
$ws.Range("E16").Value = '  -3.41%  '
$ws.Range("E17").Value = '  -2.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.029.54'
$ws.Range("E18").Value = '  -2.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.66'
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '473.76'
$ws.Range("E20").Value = '  -2.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.28'
$ws.Range("E21").Value = '  -1.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.675'
$ws.Range("E22").Value = '  -4.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.09'
$ws.Range("E23").Value = '  -1.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.25'
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.12'
$ws.Range("E25").Value = '  -2.14%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.73'
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.81'
$ws.Range("E28").Value = '  -5.28%  '
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.91'
$ws.Range("E30").Value = '  -1.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.78'
$ws.Range("E31").Value = '  -2.95%  '
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.31'
$ws.Range("E33").Value = '  -3.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '55.49'
$ws.Range("E34").Value = '  -3.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.46'
$ws.Range("E35").Value = '  +1.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.91'
$ws.Range("E36").Value = '  -2.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '458.39'
$ws.Range("E37").Value = '  -9.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.216.33'
$ws.Range("E38").Value = '  -1.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0798'
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0384'
$ws.Range("E40").Value = '  -5.90%  '
$ws.Range("E41").Value = '  -0.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.17'
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.45'
$ws.Range("E43").Value = '  -10.10%  '
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.99'
$ws.Range("E45").Value = '  +3.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.245'
$ws.Range("E46").Value = '  -5.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.00'
$ws.Range("E47").Value = '  -3.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.108'
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '118.33'
$ws.Range("E49").Value = '  -4.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₃0496'
$ws.Range("E50").Value = '  -7.72%  '
$ws.Range("E51").Value = '  +7.10%  '
